$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Clear-Cell($ws, $addr) {
    $ws.Range($addr).ClearContents()
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

$updates = @(
    ,@("H41", 11367035)
    ,@("I41", 20836306)
    ,@("K41", 20836306)
    ,@("M41", -20835866)
    ,@("H62", 2337.6)
    ,@("I62", 1513.4166)
    ,@("J62", 3573.875)
    ,@("K62", 1513.4166)
    ,@("L62", 3573.875)
    ,@("M62", -889.4166)
    ,@("N62", -4821.875)
    ,@("H65", 2337.6)
    ,@("I65", 1513.4166)
    ,@("J65", 3573.875)
    ,@("K65", 7567.083000000001)
    ,@("L65", 17869.375)
    ,@("M65", -4447.083000000001)
    ,@("N65", -24109.375)
    ,@("H74", 51736428)
    ,@("I74", 78956680)
    ,@("K74", 78956680)
    ,@("M74", -78955744)
    ,@("H77", 51736428)
    ,@("I77", 78956680)
    ,@("K77", 394783400)
    ,@("M77", -394778720)
    ,@("H97", 16675566)
    ,@("J97", 20010620)
    ,@("L97", 60031860)
    ,@("N97", -60032852)
    ,@("H98", 4868.5405)
    ,@("I98", 2861.0571)
    ,@("K98", 2861.0571)
    ,@("M98", -1363.0571)
    ,@("H107", 40279100)
    ,@("I107", 7813989)
    ,@("J107", 300000000)
    ,@("K107", 7813989)
    ,@("L107", 300000000)
    ,@("M107", -7812069)
    ,@("N107", -300003840)
    ,@("H122", 4868.5405)
    ,@("I122", 2861.0571)
    ,@("K122", 8583.1713)
    ,@("M122", -6133.1713)
    ,@("H127", 1012.5)
    ,@("I127", 985.7143)
    ,@("J127", 1200)
    ,@("K127", 2957.1429)
    ,@("L127", 3600)
    ,@("M127", 2002.8571)
    ,@("N127", -13520)
    ,@("H137", 3073.36)
    ,@("I137", 3004.6785)
    ,@("K137", 9014.0355)
    ,@("M137", -6464.0355)
    ,@("H138", 2518.61)
    ,@("I138", 2090.4211)
    ,@("J138", 2781.0483)
    ,@("K138", 6271.263300000001)
    ,@("L138", 8343.144899999999)
    ,@("M138", -1131.263300000001)
    ,@("N138", -18623.1449)
)
foreach ($u in $updates) {
    Set-Cell $ws $u[0] $u[1]
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

$updates = @(
    ,@("H74", 31250.428)
    ,@("I74", 43032.207)
    ,@("J74", 5544.727)
    ,@("K74", 43032.207)
    ,@("L74", 5544.727)
    ,@("M74", -42158.207)
    ,@("N74", -7292.727)
    ,@("H77", 31250.428)
    ,@("I77", 43032.207)
    ,@("J77", 5544.727)
    ,@("K77", 215161.035)
    ,@("L77", 27723.635)
    ,@("M77", -210793.035)
    ,@("N77", -36459.63499999999)
    ,@("H97", 4904237)
    ,@("I97", 1472.9166)
    ,@("J97", 16670871)
    ,@("K97", 1472.9166)
    ,@("L97", 16670871)
    ,@("M97", -976.9166)
    ,@("N97", -16671863)
    ,@("H102", 6063658)
    ,@("I102", 6669224)
    ,@("K102", 6669224)
    ,@("M102", -6667602)
    ,@("H132", 7080.3076)
    ,@("I132", 2587.25)
    ,@("K132", 7761.75)
    ,@("M132", -5231.75)
    ,@("H140", 53481.668)
    ,@("I140", 0)
    ,@("J140", 53481.668)
    ,@("K140", 0)
    ,@("L140", 53481.668)
    ,@("N140", -63841.668)
)
foreach ($u in $updates) {
    Set-Cell $ws $u[0] $u[1]
}
$clears = @(
    "M140"
)
foreach ($addr in $clears) {
    Clear-Cell $ws $addr
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

$updates = @(
    ,@("H86", 8099539)
    ,@("I86", 10915056)
    ,@("K86", 10915056)
    ,@("M86", -10913933)
    ,@("H89", 8099539)
    ,@("I89", 10915056)
    ,@("K89", 54575280)
    ,@("M89", -54569664)
    ,@("H94", 2512.303)
    ,@("I94", 1183.174)
    ,@("K94", 1183.174)
    ,@("M94", -732.174)
    ,@("H134", 7817004.5)
    ,@("I134", 16667946)
    ,@("K134", 50003838)
    ,@("M134", -50001303)
)
foreach ($u in $updates) {
    Set-Cell $ws $u[0] $u[1]
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

$updates = @(
    ,@("H48", 0)
    ,@("J48", 0)
    ,@("L48", 0)
    ,@("H58", 3826.1746)
    ,@("I58", 2562.4167)
    ,@("J58", 5511.185)
    ,@("K58", 2562.4167)
    ,@("L58", 5511.185)
    ,@("M58", -2359.4167)
    ,@("N58", -5917.185)
    ,@("H136", 3826.1746)
    ,@("I136", 2562.4167)
    ,@("J136", 5511.185)
    ,@("K136", 7687.250100000001)
    ,@("L136", 16533.555)
    ,@("M136", -5137.250100000001)
    ,@("N136", -21633.555)
)
foreach ($u in $updates) {
    Set-Cell $ws $u[0] $u[1]
}
$clears = @(
    "N48"
)
foreach ($addr in $clears) {
    Clear-Cell $ws $addr
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

$updates = @(
    ,@("H5", 1019.67645)
    ,@("I5", 696.5789)
    ,@("J5", 1428.9333)
    ,@("K5", 2089.7367)
    ,@("L5", 4286.7999)
    ,@("M5", -1977.7367)
    ,@("N5", -4510.7999)
    ,@("H32", 158.16667)
    ,@("I32", 0)
    ,@("J32", 158.16667)
    ,@("K32", 0)
    ,@("L32", 474.50001)
    ,@("N32", -1040.50001)
    ,@("H55", 7707446.5)
    ,@("J55", 10016800)
    ,@("L55", 30050400)
    ,@("N55", -30050754)
    ,@("H113", 8923.916999999999)
    ,@("J113", 13942)
    ,@("L113", 41826)
    ,@("N113", -46166)
    ,@("H135", 1019.67645)
    ,@("I135", 696.5789)
    ,@("J135", 1428.9333)
    ,@("K135", 6269.2101)
    ,@("L135", 12860.3997)
    ,@("M135", -3734.2101)
    ,@("N135", -17930.3997)
)
foreach ($u in $updates) {
    Set-Cell $ws $u[0] $u[1]
}
$clears = @(
    "M32"
)
foreach ($addr in $clears) {
    Clear-Cell $ws $addr
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

$updates = @(
    ,@("H70", 6085.404)
    ,@("I70", 5550.2285)
    ,@("K70", 5550.2285)
    ,@("M70", -5280.2285)
    ,@("H73", 6085.404)
    ,@("I73", 5550.2285)
    ,@("K73", 5550.2285)
    ,@("M73", -4614.2285)
    ,@("H93", 59993.332)
    ,@("J93", 59993.332)
    ,@("L93", 59993.332)
    ,@("N93", -63737.332)
    ,@("H97", 1771.5555)
    ,@("I97", 920)
    ,@("K97", 920)
    ,@("M97", -424)
    ,@("H102", 3831.8)
    ,@("I102", 3831.8)
    ,@("K102", 3831.8)
    ,@("M102", -2209.8)
    ,@("H109", 60284)
    ,@("J109", 60284)
    ,@("L109", 60284)
    ,@("N109", -62364)
    ,@("H132", 3738.25)
    ,@("I132", 2985.1667)
    ,@("J132", 5997.5)
    ,@("K132", 8955.500100000001)
    ,@("L132", 17992.5)
    ,@("M132", -6425.500100000001)
    ,@("N132", -23052.5)
)
foreach ($u in $updates) {
    Set-Cell $ws $u[0] $u[1]
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

$updates = @(
    ,@("H46", 2677.7407)
    ,@("I46", 2007.0714)
    ,@("K46", 2007.0714)
    ,@("M46", -1819.0714)
    ,@("H93", 1041.1428)
    ,@("I93", 1131.4546)
    ,@("K93", 1131.4546)
    ,@("M93", 116.5454)
    ,@("H100", 5117.6)
    ,@("J100", 6399.5)
    ,@("L100", 6399.5)
    ,@("N100", -7481.5)
    ,@("H132", 10642916)
    ,@("I132", 20003080)
    ,@("J132", 6366.227)
    ,@("K132", 60009240)
    ,@("L132", 19098.681)
    ,@("M132", -60006710)
    ,@("N132", -24158.681)
    ,@("H136", 8075.0864)
    ,@("I136", 2780.8438)
    ,@("J136", 14591.077)
    ,@("K136", 8342.5314)
    ,@("L136", 43773.231)
    ,@("M136", -5792.5314)
    ,@("N136", -48873.231)
)
foreach ($u in $updates) {
    Set-Cell $ws $u[0] $u[1]
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

$updates = @(
    ,@("H97", 0)
    ,@("J97", 0)
    ,@("L97", 0)
    ,@("H132", 8099.1304)
    ,@("I132", 10749.182)
    ,@("J132", 5669.9165)
    ,@("K132", 32247.546)
    ,@("L132", 17009.7495)
    ,@("M132", -29717.546)
    ,@("N132", -22069.7495)
)
foreach ($u in $updates) {
    Set-Cell $ws $u[0] $u[1]
}
$clears = @(
    "N97"
)
foreach ($addr in $clears) {
    Clear-Cell $ws $addr
}
